$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "39.445.19"
$ws.Range("E2").Value = "  -2.91%  "
Set-TextValue $ws.Range("D3") "2.219.10"
$ws.Range("E3").Value = "  -6.28%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "297.03"
$ws.Range("E5").Value = "  -4.39%  "
Set-TextValue $ws.Range("D6") "83.05"
$ws.Range("E6").Value = "  -3.70%  "
Set-TextValue $ws.Range("D7") "0.511"
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.469"
$ws.Range("E9").Value = "  -4.36%  "
Set-TextValue $ws.Range("D10") "0.0775"
$ws.Range("E10").Value = "  -7.69%  "
Set-TextValue $ws.Range("D11") "29.16"
$ws.Range("E11").Value = "  -3.56%  "
Set-TextValue $ws.Range("D12") "47.83"
$ws.Range("E12").Value = "  -9.34%  "
$ws.Range("E13").Value = "  -2.08%  "
Set-TextValue $ws.Range("D14") "2.561.21"
$ws.Range("E14").Value = "  -6.41%  "
Set-TextValue $ws.Range("D15") "6.31"
$ws.Range("E15").Value = "  -3.15%  "
Set-TextValue $ws.Range("D16") "14.11"
$ws.Range("E16").Value = "  -5.44%  "
Set-TextValue $ws.Range("D17") "2.206.86"
$ws.Range("E17").Value = "  -7.93%  "
Set-TextValue $ws.Range("D18") "0.716"
$ws.Range("E18").Value = "  -5.19%  "
Set-TextValue $ws.Range("D19") "39.341.46"
$ws.Range("E19").Value = "  -3.08%  "
Set-TextValue $ws.Range("D20") "0.0₃0874"
$ws.Range("E20").Value = "  -3.85%  "
Set-TextValue $ws.Range("D21") "5.72"
$ws.Range("E21").Value = "  -6.38%  "
Set-TextValue $ws.Range("D22") "65.06"
$ws.Range("E22").Value = "  -4.67%  "
Set-TextValue $ws.Range("D23") "10.27"
$ws.Range("E23").Value = "  -4.22%  "
Set-TextValue $ws.Range("D24") "227.77"
$ws.Range("E24").Value = "  -3.01%  "
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  -0.11%  "
Set-TextValue $ws.Range("D26") "2.41"
$ws.Range("E26").Value = "  -6.48%  "
Set-TextValue $ws.Range("D27") "1.81"
$ws.Range("E27").Value = "  +0.96%  "
Set-TextValue $ws.Range("D28") "22.58"
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("E29").Value = "  -2.50%  "
Set-TextValue $ws.Range("D30") "9.13"
$ws.Range("E30").Value = "  -0.88%  "
Set-TextValue $ws.Range("D31") "149.63"
$ws.Range("E31").Value = "  -2.58%  "
Set-TextValue $ws.Range("D32") "32.13"
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("E33").Value = "  -0.13%  "
Set-TextValue $ws.Range("D34") "4.85"
$ws.Range("E34").Value = "  -6.11%  "
Set-TextValue $ws.Range("D35") "0.0694"
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("E36").Value = "  -3.21%  "
Set-TextValue $ws.Range("D37") "0.110"
$ws.Range("E37").Value = "  -3.49%  "
Set-TextValue $ws.Range("D38") "0.0966"
$ws.Range("E38").Value = "  -3.36%  "
Set-TextValue $ws.Range("D39") "15.26"
$ws.Range("E39").Value = "  -4.21%  "
Set-TextValue $ws.Range("D40") "2.64"
$ws.Range("E40").Value = "  -4.07%  "
$ws.Range("E41").Value = "  -2.61%  "
Set-TextValue $ws.Range("D42") "3.66"
$ws.Range("E42").Value = "  -4.33%  "
Set-TextValue $ws.Range("D43") "1.908.84"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("E44").Value = "  -3.67%  "
Set-TextValue $ws.Range("D45") "2.03"
$ws.Range("E45").Value = "  -15.27%  "
Set-TextValue $ws.Range("D46") "16.07"
$ws.Range("E46").Value = "  -8.79%  "
Set-TextValue $ws.Range("D47") "8.97"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("E48").Value = "  -2.32%  "
Set-TextValue $ws.Range("D49") "2.431.69"
$ws.Range("E49").Value = "  -6.47%  "
Set-TextValue $ws.Range("D50") "70.66"
$ws.Range("E50").Value = "  -1.31%  "
Set-TextValue $ws.Range("D51") "87.07"
$ws.Range("E51").Value = "  -6.19%  "
